$d = $word.ActiveDocument

$map = [ordered]@{
    "73×86=" = "16×23="
    "49×97=" = "85×85="
    "84×92=" = "41×62="
    "76×36=" = "95×51="
    "19×35=" = "34×84="
    "42×47=" = "51×55="
    "74×53=" = "40×51="
    "71×63=" = "48×17="
    "99×95=" = "97×85="
    "53×30=" = "19×68="
    "51×93=" = "73×66="
    "40×12=" = "41×97="
    "94×45=" = "20×94="
    "67×85=" = "93×60="
    "39×66=" = "53×52="
    "46×42=" = "41×20="
    "53×32=" = "46×18="
    "11×91=" = "74×78="
    "43×25=" = "49×68="
    "36×63=" = "74×78="
    "56×76=" = "90×23="
    "40×46=" = "52×57="
    "72×91=" = "84×65="
    "72×49=" = "99×85="
    "11×97=" = "11×76="
}

foreach ($key in $map.Keys) {
    $d.Content.Find.Execute($key, $true, $false, $false, $false, $false, $true, 1, $false, $map[$key], 2)
}
